$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the last existing row (232) down onto the new rows (233:244)
$ws.Range("A232:I232").Copy()
$ws.Range("A233:I244").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A233").Value = 45898
$ws.Range("B233").Value = 'Levy Ndoutoume'
$ws.Range("C233").Value = 55
$ws.Range("D233").Value = 6
$ws.Range("E233").Value = 6
$ws.Range("F233").Value = 4
$ws.Range("G233").Value = 'Ischio'
$ws.Range("H233").Value = 5
$ws.Range("I233").Formula = "=C233*D233"

$ws.Range("A234").Value = 45898
$ws.Range("B234").Value = 'Jeremie Laurent'
$ws.Range("C234").Value = 55
$ws.Range("D234").Value = 6
$ws.Range("E234").Value = 2
$ws.Range("F234").Value = 0
$ws.Range("G234").Value = ""
$ws.Range("H234").Value = 8
$ws.Range("I234").Formula = "=C234*D234"

$ws.Range("A235").Value = 45898
$ws.Range("B235").Value = 'Naim Dhib'
$ws.Range("C235").Value = 55
$ws.Range("D235").Value = 4
$ws.Range("E235").Value = 3
$ws.Range("F235").Value = 1
$ws.Range("G235").Value = 'Courbature '
$ws.Range("H235").Value = 7
$ws.Range("I235").Formula = "=C235*D235"

$ws.Range("A236").Value = 45898
$ws.Range("B236").Value = 'Amine Taiar'
$ws.Range("C236").Value = 55
$ws.Range("D236").Value = 3
$ws.Range("E236").Value = 5
$ws.Range("F236").Value = 6
$ws.Range("G236").Value = 'Dos'
$ws.Range("H236").Value = 7
$ws.Range("I236").Formula = "=C236*D236"

$ws.Range("A237").Value = 45898
$ws.Range("B237").Value = 'Mattheo Haon'
$ws.Range("C237").Value = 55
$ws.Range("D237").Value = 4
$ws.Range("E237").Value = 5
$ws.Range("F237").Value = 0
$ws.Range("G237").Value = ""
$ws.Range("H237").Value = 6
$ws.Range("I237").Formula = "=C237*D237"

$ws.Range("A238").Value = 45898
$ws.Range("B238").Value = 'Karahali Souaré'
$ws.Range("C238").Value = 55
$ws.Range("D238").Value = 3
$ws.Range("E238").Value = 5
$ws.Range("F238").Value = 7
$ws.Range("G238").Value = 'Cheville'
$ws.Range("H238").Value = 3
$ws.Range("I238").Formula = "=C238*D238"

$ws.Range("A239").Value = 45898
$ws.Range("B239").Value = 'Emmanuel Valey'
$ws.Range("C239").Value = 55
$ws.Range("D239").Value = 3
$ws.Range("E239").Value = 4
$ws.Range("F239").Value = 4
$ws.Range("G239").Value = 'Adducteur'
$ws.Range("H239").Value = 6
$ws.Range("I239").Formula = "=C239*D239"

$ws.Range("A240").Value = 45898
$ws.Range("B240").Value = 'Naim Ighbane'
$ws.Range("C240").Value = 55
$ws.Range("D240").Value = 3
$ws.Range("E240").Value = 3
$ws.Range("F240").Value = 0
$ws.Range("G240").Value = ""
$ws.Range("H240").Value = 7
$ws.Range("I240").Formula = "=C240*D240"

$ws.Range("A241").Value = 45898
$ws.Range("B241").Value = 'Omar Benyounes'
$ws.Range("C241").Value = 55
$ws.Range("D241").Value = 5
$ws.Range("E241").Value = 5
$ws.Range("F241").Value = 0
$ws.Range("G241").Value = ""
$ws.Range("H241").Value = 7
$ws.Range("I241").Formula = "=C241*D241"

$ws.Range("A242").Value = 45898
$ws.Range("B242").Value = 'Romain Thunet'
$ws.Range("C242").Value = 55
$ws.Range("D242").Value = 6
$ws.Range("E242").Value = 5
$ws.Range("F242").Value = 2
$ws.Range("G242").Value = 'Orteil'
$ws.Range("H242").Value = 6
$ws.Range("I242").Formula = "=C242*D242"

$ws.Range("A243").Value = 45898
$ws.Range("B243").Value = 'Yanis Berrached'
$ws.Range("C243").Value = 55
$ws.Range("D243").Value = 5
$ws.Range("E243").Value = 5
$ws.Range("F243").Value = 0
$ws.Range("G243").Value = ""
$ws.Range("H243").Value = 3
$ws.Range("I243").Formula = "=C243*D243"

$ws.Range("A244").Value = 45898
$ws.Range("B244").Value = 'Amir Etien'
$ws.Range("C244").Value = 55
$ws.Range("D244").Value = 5
$ws.Range("E244").Value = 6
$ws.Range("F244").Value = 4
$ws.Range("G244").Value = 'Quadri'
$ws.Range("H244").Value = 0
$ws.Range("I244").Formula = "=C244*D244"

# Move the visible selection to match the updated view (K240) and scroll the window
$ws.Range("K240").Select()
$excel.ActiveWindow.ScrollRow = 216